# Apply the "AquiferOpenStudyNotesBookIntros" resource-data update.
#
# Summary of the edit:
#  1. The "License Information" Heading2 paragraph loses its heading
#     style and its text is replaced by the bold title
#     "Aquifer Open Study Notes (Book Intros)".
#  2. The license paragraph right after it is rewritten: the old
#     "Notes d'etude..." bold lead-in + "(French) is based on: " is
#     replaced by a plain "This work is an adaptation of " lead-in,
#     "Tyndale Open Study Notes" stays, and the old
#     ", Tyndale House Publishers [hyperlink], 2019, which is licensed
#     under a CC BY-SA 4.0 license [hyperlink]." tail (which included
#     two hyperlinks) is replaced by new plain-text license wording
#     that no longer links out (the two hyperlinks are removed).
#  3. The "This PDF version is provided under the same license."
#     paragraph is replaced with the multi-language availability blurb.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. "License Information" heading -> plain bold title paragraph
# ---------------------------------------------------------------
$pHeading = $d.Paragraphs.Item(4)
$pHeading.Range.Find.Execute("License Information", $false, $false, $false, `
    $false, $false, $true, 1, $false, "Aquifer Open Study Notes (Book Intros)", 2)
$pHeading.Range.ParagraphFormat.Style = $d.Styles.Item("Normal")
$titleRng = $d.Range($pHeading.Range.Start, $pHeading.Range.End - 1)
$titleRng.Font.Bold = $true

# ---------------------------------------------------------------
# 2. License paragraph rewrite
# ---------------------------------------------------------------
$pLicense = $d.Paragraphs.Item(5)

# Remove the two hyperlinks that live in this paragraph (Tyndale House
# Publishers / CC BY-SA 4.0 license). Hyperlink.Delete() strips the
# hyperlink field but keeps the visible text as a plain run, which we
# then fold into the text replacement below.
$d.Hyperlinks.Item(1).Delete()
$d.Hyperlinks.Item(1).Delete()

# Replace the bold lead-in + "(French) is based on: " with the new
# plain lead-in sentence, then strip the bold formatting it inherited.
$pLicense.Range.Find.Execute( `
    "Notes d'étude - Introductions aux livres (Tyndale) (French) is based on: ", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "This work is an adaptation of ", 2)

$boldRng = $pLicense.Range.Duplicate
$boldRng.Find.Execute("This work is an adaptation of ")
$boldRng.Font.Bold = $false

# Replace the remainder (the now-plain former hyperlink text, the
# publisher/year/license sentence and the final period) with the new
# attribution text.
$pLicense.Range.Find.Execute( `
    ", Tyndale House Publishers, 2019, which is licensed under a CC BY-SA 4.0 license.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    " © 2023 Tyndale House Publishers, licensed under the CC BY-SA 4.0 license. The adaptation, Aquifer Open Study Notes, was created by Mission Mutual and is also licensed under CC BY-SA 4.0.", 2)

# ---------------------------------------------------------------
# 3. "This PDF version..." -> multi-language availability blurb
# ---------------------------------------------------------------
$pPdf = $d.Paragraphs.Item(6)
$pPdf.Range.Find.Execute( `
    "This PDF version is provided under the same license.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "This resource has been adapted into multiple languages, including English, Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文).", 2)
